$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at row 227 (Femacal de La Calera,
# Apio, Primera quality, week of 2021-11-05). All existing records from the
# old row 227 onward shift down by one row (227->228, ..., 262->263), which
# $ws.Rows.Item(227).Insert() handles for us.
$ws.Rows.Item(227).Insert()

$ws.Cells.Item(227, 1).Value = 3
$ws.Cells.Item(227, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(227, 3).Value = 'Coquimbo'
$ws.Cells.Item(227, 4).Value = 44505
$ws.Cells.Item(227, 5).Value = 5
$ws.Cells.Item(227, 6).Value = 100112017
$ws.Cells.Item(227, 7).Value = 'Apio'
$ws.Cells.Item(227, 8).Value = 'Americana (o)'
$ws.Cells.Item(227, 9).Value = 'Primera'
$ws.Cells.Item(227, 10).Value = 160
$ws.Cells.Item(227, 11).Value = 9000
$ws.Cells.Item(227, 12).Value = 9000
$ws.Cells.Item(227, 13).Value = 9000
$ws.Cells.Item(227, 14).Value = '$/docena de matas'
$ws.Cells.Item(227, 15).Value = 'Pan de Azúcar'
$ws.Cells.Item(227, 16).Value = 1500
$ws.Cells.Item(227, 17).Value = 6
$ws.Cells.Item(227, 18).Value = 'Hortaliza'
